# "1st changes of mifos to finflux"
# The "Repayment schedule" sheet gets a new (blank) column inserted right
# before column N ("Late"/"Outstanding" group), shifting the existing
# N:P ("Late", "heading"/Original, "Outstanding") block one column to the
# right (into O:Q). The sheet also becomes the active tab/selection instead
# of "Transactions".

$wb = $excel.ActiveWorkbook

# --- Repayment schedule: insert a blank column before column N -----------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("N").Insert() | Out-Null

# New column keeps roughly the same width as its left neighbour (column M).
$wsSchedule.Columns("N").ColumnWidth = 9.7

# --- Make "Repayment schedule" the active sheet/selection -----------------
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("R5").Select() | Out-Null
